$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet and update its "through" date label
$ws.Name = "Through 2022-10-05"
$ws.Range("I1").Value = "2022 (through 10-05)"

# Update the October figures (row 11) and Total row (row 14) for column I
$ws.Range("I11").Value = 17
$ws.Range("I14").Value = 1299
